# Generate Report for Handback
# Updates status/timestamp cells across the Overview, zh-cn and de-de sheets
# to reflect a newer handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
# Rows 3 and 5 both show the 6d657a99 / e825220c timestamp that gets bumped.
$wsOverview.Range("G3").Value = "2016-08-24 22:16:58"
$wsOverview.Range("G5").Value = "2016-08-24 22:16:58"

# --- zh-cn sheet ---
# Priority column (E): "ht" -> "mt" for the affected rows
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-08-24 22:16:53"
$wsZhCn.Range("K3").Value = "2016-08-24 22:17:15"

# --- de-de sheet ---
# Priority column (E): "ht" -> "mt" for the affected rows
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsDeDe.Range("H3").Value = "2016-08-24 22:16:58"
$wsDeDe.Range("K3").Value = "2016-08-24 22:17:22"
